# Atualiza Caixa 25 com roupas do Bento
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Catálogo de Mudança")

# Insert a new row at position 196 (pushes existing row 196.. down to 197..)
$ws.Rows.Item(196).Insert()

# Copy the formatting (styles/fills/borders) from the row that is now 197
# (the former row 196) down into the newly blank row 196, so the new row
# keeps the same visual styling as the rest of the table.
$ws.Range("A197:F197").Copy()
$ws.Range("A196:F196").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new row with the new catalog entry
$ws.Range("A196").Value = "Caixa 25"
$ws.Range("B196").Value = "Vestuário infantil"
$ws.Range("C196").Value = "Roupas diversas do Bento (lavadas recentemente)"
$ws.Range("D196").Value = "Quarto do Bento"
$ws.Range("E196").Value = "Alta"
$ws.Range("F196").Value = ""

# Update the summary sheet counts
$resumo = $wb.Worksheets.Item("Resumo")
$resumo.Range("B3").Value = 437
$resumo.Range("B6").Value = 96
